# TC_30.xlsx edit: rename sheet, refresh embedded CEIC comment payload,
# widen the decimal number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "My Series" to "Data".
$ws.Name = "Data"

# 2. Replace the CEIC add-in metadata blob stored in the A1 cell comment.
$comment = $ws.Range("A1").Comment
$comment.Text("GRoAAB+LCAAAAAAAAAOlGdtu48b1Vwg9tUAlUpLXsY1ZBrrZESpZhiTX67wUI3JsTU1xFHJoWW8pkCBFmqIoik2RXtGnFAW6XbQJkO728i/B2rt96i/0zIVDUpKz1XZhrDnnNmfOnOsYvX0zC6xrEsWUhQ9L1YpTskjoMZ+Glw9LCb8oV3dLb7uoc+OR4ARHeEY4EFvAFcYHNzF9WJpyPj+w7cViUVnUKyy6tGuOU7Uf9Xsjb0pmuEzDmOPQIyXD5b+eq+Silj/rE459zLHifFjqjrqVFqFeG2B9HOJLElWaSUxDEsedkFNOSSw4I4I5abX731MHc2uV3UoV2WvwjLKZ0MBXdAVKBdd0sC0Z0xlxa051r+zslevVsbN74NQO6k6lVtt5N2U0hKiHYz4i0TX1JGDE8Wwu2Z29etXZdWp1B9kbiUBWZgAXDQJ/SK5pTPwWCYJ4K4vY+gIbHodTb2dMUC/HqwW9uQpHEZ5Px5QHZDs1hv2mNQu1LpkQFx2yiHhgvzdS6ZgsBpE263jeA+x4SiO+bOPl1rJOYxIN5sJI27G6qM1C3ghIxE/ncNfEB1cAhMujhCD7HmTG1KaxB980TIjvXuAgzjMVkOiMRVfxHHvkGOLYFjIWYcCwDw7HacypF2cC1jDoJGJzEAm7N1ngH4LYlHoDxsjuhmBksXGTsatV8UUkkvcqbxhudYaN/DU4Gk3ZYhAGy1Eyib2ITojfbqbUG3FIhKTmbiUxZzPQIgMhBctBlvAPQnAVjNrEozMcnARgxtitg5QCADUSzi4ob7EgmYXGnitQdAYnGpMbc0KzRgO43lBYnYXdMKVXrrARVWQYsoXZcx0hjZADN2Iv9bJ1xCpxG2Dp9a1j5I2IUx7SACpE/i5y0KJXjKaE8I0uoTBIJMNDUXPc5vI4mU0gwCYQZddy1xjZGR6Bo4K3g16uA3WkLH/GjnMgf0APg0ad0L+fLkUi2C63l1sF3AoIwZmCZoDDK4CeUT49bqRn2YBBygL30q/jEMTuPMBLCTZWysNQN/SCxCcqJXTDC+miQjdNfS8erYF6EOUuwuFyvJyL1GDfQ9EmFzgJoFpxSCqXWb5YAaNGfLVKkweh0yhIb9wVvUAMzYDnzyoepEtR8CoemwmADTX4bITsPL3I+R7phJc9HF4mkFWNHVfhxt9EPhhHOIzFcUwSXXG9zUQovReV3F11WYNEupO6LAZYZK/QoTGZzVmEgz4Yhh4moSygukKAhfuYT/UKYjkgXmpkO2M1XEXNUsVfRyaDUh1DZPx8UOagkkocRvUdKiEUYUgcs898qPw4oJOo4GUbcXBjWT1M/U2cbsvamF4C9J6Qbr5LlqIhyRYaLpzWraYIuUItUMkdDXf2ag+ceg1SuVgjeeIhwYHVga6UE6sbXpOYz4DtwBqSmPrwRXFwYL1DJoRC1EsT6bjbmjvPhw4j8l4C3fVSqtIAfYuQIgEk2Esa4mCd0GAyBvec4ChY5gjVUXvMA7q7H/3r9tdPXzz7/O7jx6++/OA/f//li3/89PbJh/Bx95e/3n7yC3VMRYzGeBIQqdC4ubfn1HfA0QwI6fwAJd9PPC5h5+eyBTBrpHtZuWh1uq2jXlMmFANM2UV7kYhg6OElS7LlSB1CbiSv1E49QZG44zRB6XUB2yai7qsoFG3rNSlS5/H3MSpbvHz++cvnf7qXWxssKy7V/f0H5WrttbUHBoDqGp2pPWIcUN1eOmzslJ0H5VotR7xCg4YwM0FzZ+zU9V2YL/bFfFE1ydw3jryJaBWlJY3xpb3Cp0AtloQ8WhoXyK9TpHT8MYSIQatQyC20i37541d/flyg0tbVkKIUUI4lkfImO11I0cfDsTUanA5bHWvcGQk/yXA5OiX8G4j17iaeCk4VhgkOvgOTspiTrRK0gCWLXVgEe1NrCZGYi8OCs22Cqo3eUOSqlkcRS+bqRnIMGXQDpckmGzk25BqJk/ZcSzoZagO50vX2b19sYtAH0W52GlJuRr48DBUwCpTD66j97J8vvvroxbNnd09/dvvVDwsS9D5mEgI/h2jKL43bQ8rT9WYFgs5G0phXzvdz9UUDRed4wmjIY7e6K5tGvULAWhXS5G/UnUHJk4KlvQC+AkHv4Lhzw3Vgu8fILgJAzzmGasuyZtsAVA7P7Prv3/z27ldf3H369NVHf7z9+A+3n3z68vnvXj35vYq6u8dP737yRGf51UIgdREtvOoCLTkRepaIRkvUbuvr939uhYxb0HNYicxIX7//WU6YUFR2J5lk6OmMIkUV1kjzzILPyqlidCjwGRbVALRECasbCl3E2Jx62SbvloUoEXcS8a3uuJzExGLQTn0bTlIkzpj/Vz7NokrqyVtOrVrTWKWNOMIExznTHwVsAk1GipAT1wpJgeubGTJaud9Rb9Bs9DISpcQg8mFAc8TYKD5Q2lOKktKN05WZLDIIYKHx85JADMlrZOsoIzmXxmw9bV40fJH+Ns9mBQoY2KNINUShfr4cJXPohrluYO/Hy5eaXP97rHrVfEecrbvtIh7WOSwUwiJaACRepiaNUmmqG4vBVrWzx8I02RJwhdcdMId+oVSd1jX0lZEt8k4nili0MflkmJSsD500ZBQ7s7ihkXequm4/u6sUkCY886FGP31C1iYB4ds939kZd59dvzEv3P22rN14EPjamNuNHsYsmYD8G6ZwlP/3CVM5WyOKoLESLx5bvzmmk+sQBt4ttVFHkYxiAoTd9aPgIY1i/khkAv2lIOcGcq461Edi4FIfcn3u7qs14O28cLugZRq5XD0ks6BHZ3TLqdBJw7soBEw5n6sOrrudo4jKckxuoL/MSYCcOPkBVA0x8WwnTfkrpFLDL95qYno55dsq9tYEE59MnLI3IbXyju/slfcJqZerVfgfe7Wa4zwQDz1aOCQOShZbbmKnF5b9gcf9Lx+K5dsZGgAA")

# 3. Widen the custom decimal number format used by column B
#    (0.000 -> ###0.000).
$ws.Range("B1:B10").NumberFormat = "###0.000"
